$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2028.125
$ws.Range("I28").Value = 1351.7778
$ws.Range("J28").Value = 2897.7144
$ws.Range("K28").Value = 1351.7778
$ws.Range("L28").Value = 2897.7144
$ws.Range("M28").Value = -866.7778000000001
$ws.Range("N28").Value = -3867.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 19285.715
$ws.Range("J23").Value = 19285.715
$ws.Range("L23").Value = 19285.715
$ws.Range("N23").Value = -19803.715

$ws.Range("H32").Value = 1395.29
$ws.Range("I32").Value = 1381.101
$ws.Range("J32").Value = 2800
$ws.Range("K32").Value = 1381.101
$ws.Range("L32").Value = 2800
$ws.Range("M32").Value = -1094.101
$ws.Range("N32").Value = -3374

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 851.25
$ws.Range("I8").Value = 700
$ws.Range("J8").Value = 1002.5
$ws.Range("K8").Value = 700
$ws.Range("L8").Value = 1002.5
$ws.Range("M8").Value = -560
$ws.Range("N8").Value = -1282.5

$ws.Range("H22").Value = 677.5
$ws.Range("I22").Value = 658.8889
$ws.Range("J22").Value = 733.3333
$ws.Range("K22").Value = 658.8889
$ws.Range("L22").Value = 733.3333
$ws.Range("M22").Value = -485.8889
$ws.Range("N22").Value = -1079.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 821.44446
$ws.Range("I16").Value = 600
$ws.Range("J16").Value = 932.1667
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 932.1667
$ws.Range("M16").Value = -313
$ws.Range("N16").Value = -1506.1667

$ws.Range("H31").Value = 2186.1353
$ws.Range("I31").Value = 1440.2963
$ws.Range("J31").Value = 4199.9
$ws.Range("K31").Value = 1440.2963
$ws.Range("L31").Value = 4199.9
$ws.Range("M31").Value = -1145.2963
$ws.Range("N31").Value = -4789.9

$ws.Range("H34").Value = 2186.1353
$ws.Range("I34").Value = 1440.2963
$ws.Range("J34").Value = 4199.9
$ws.Range("K34").Value = 1440.2963
$ws.Range("L34").Value = 4199.9
$ws.Range("M34").Value = -1238.2963
$ws.Range("N34").Value = -4603.9

$ws.Range("H107").Value = 1077.75
$ws.Range("I107").Value = 604.8125
$ws.Range("J107").Value = 1708.3334
$ws.Range("K107").Value = 604.8125
$ws.Range("L107").Value = 1708.3334
$ws.Range("M107").Value = 1315.1875
$ws.Range("N107").Value = -5548.3334

$ws.Range("H113").Value = 821.44446
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 932.1667
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 932.1667
$ws.Range("M113").Value = 1570
$ws.Range("N113").Value = -5272.1667

$ws.Range("H132").Value = 1612327.8
$ws.Range("I132").Value = 1900.3334
$ws.Range("J132").Value = 4631879
$ws.Range("K132").Value = 5701.0002
$ws.Range("L132").Value = 13895637
$ws.Range("M132").Value = -3171.0002
$ws.Range("N132").Value = -13900697

$ws.Range("H141").Value = 63227.145
$ws.Range("J141").Value = 63227.145
$ws.Range("L141").Value = 63227.145
$ws.Range("N141").Value = -73587.14499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2758
$ws.Range("I25").Value = 800
$ws.Range("J25").Value = 3247.5
$ws.Range("K25").Value = 2400
$ws.Range("L25").Value = 9742.5
$ws.Range("M25").Value = -2231
$ws.Range("N25").Value = -10080.5

$ws.Range("H30").Value = 2758
$ws.Range("I30").Value = 800
$ws.Range("J30").Value = 3247.5
$ws.Range("K30").Value = 2400
$ws.Range("L30").Value = 9742.5
$ws.Range("M30").Value = -2298
$ws.Range("N30").Value = -9946.5

$ws.Range("H122").Value = 23939.936
$ws.Range("J122").Value = 1204.7778
$ws.Range("L122").Value = 10843.0002
$ws.Range("N122").Value = -15743.0002

$ws.Range("H123").Value = 3860.8333
$ws.Range("J123").Value = 5157.143
$ws.Range("L123").Value = 15471.429
$ws.Range("N123").Value = -20371.429

$ws.Range("H131").Value = 1326.6737
$ws.Range("J131").Value = 1258.4945
$ws.Range("L131").Value = 3775.4835
$ws.Range("N131").Value = -13855.4835

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 31250
$ws.Range("I4").Value = 2500
$ws.Range("K4").Value = 2500
$ws.Range("M4").Value = -2388

$ws.Range("H20").Value = 5004500
$ws.Range("I20").Value = 10000000
$ws.Range("K20").Value = 10000000
$ws.Range("M20").Value = -9999755

$ws.Range("H122").Value = 951.1111
$ws.Range("I122").Value = 872
$ws.Range("J122").Value = 1050
$ws.Range("K122").Value = 2616
$ws.Range("L122").Value = 3150
$ws.Range("M122").Value = -166
$ws.Range("N122").Value = -8050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 42437.043
$ws.Range("I22").Value = 143282.72
$ws.Range("J22").Value = 912.35297
$ws.Range("K22").Value = 143282.72
$ws.Range("L22").Value = 912.35297
$ws.Range("M22").Value = -142987.72
$ws.Range("N22").Value = -1502.35297

$ws.Range("H27").Value = 42437.043
$ws.Range("I27").Value = 143282.72
$ws.Range("J27").Value = 912.35297
$ws.Range("K27").Value = 143282.72
$ws.Range("L27").Value = 912.35297
$ws.Range("M27").Value = -143175.72
$ws.Range("N27").Value = -1126.35297

$ws.Range("H33").Value = 4899
$ws.Range("I33").Value = 4899
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 4899
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -4609
$ws.Range("N33").ClearContents()

$ws.Range("H122").Value = 2708.1875
$ws.Range("I122").Value = 2128.25
$ws.Range("J122").Value = 3288.125
$ws.Range("K122").Value = 6384.75
$ws.Range("L122").Value = 9864.375
$ws.Range("M122").Value = -3934.75
$ws.Range("N122").Value = -14764.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 11753.75
$ws.Range("I21").Value = 8507.5
$ws.Range("K21").Value = 8507.5
$ws.Range("M21").Value = -8272.5

$ws.Range("H35").Value = 11753.75
$ws.Range("I35").Value = 8507.5
$ws.Range("K35").Value = 8507.5
$ws.Range("M35").Value = -8217.5

$ws.Range("H42").Value = 12500
$ws.Range("J42").Value = 12500
$ws.Range("L42").Value = 12500
$ws.Range("N42").Value = -13256

$ws.Range("H43").Value = 8506.6
$ws.Range("I43").Value = 8133.25
$ws.Range("K43").Value = 8133.25
$ws.Range("M43").Value = -7984.25
